$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "How do I sign into Office 365"
$ws.Range("B1").Value = "The reception phone number is 403-298-2200"
$ws.Range("A1").Value = "What is the Reception Phone Number?`nWhat is the main line`nWhat is the Enerplus phone number"
$ws.Range("B2").Value = "Go to the following site https://portal.office.com"

$ws.Range("A1").WrapText = $true

$ws.Columns.Item(1).ColumnWidth = 34
$ws.Columns.Item(2).ColumnWidth = 29

$ws.Rows.Item(1).RowHeight = 43.5

$ws.Range("B11").Select() | Out-Null
